# Update the daily scrum log: 2 story points were completed on day 2 (row 8).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Completed" value for 7/9 (row 8) from 0 to 2 story points.
$ws.Range("D8").Value = 2

# Excel recalculates dependent formulas (E8:E10, E12, B12, C12) automatically.
$excel.Calculate()

# Move the active cell selection to D14, matching where the author left off
# after logging today's scrum update.
$ws.Range("D14").Select()
